# Mise à jour des résultats du script
# Appends new scraped rows (2025-10-13) to the data table on Sheet1,
# extending it from A1:D264 to A1:D300.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 265

$data = @(
    @('2025-10-13', 'eaux souterraines', 100, 1),
    @('2025-10-13', 'zone tampon', 102, 4),
    @('2025-10-13', 'zone tampon', 103, 2),
    @('2025-10-13', 'eaux de surface', 103, 1),
    @('2025-10-13', 'eaux souterraines', 104, 4),
    @('2025-10-13', 'eaux de surface', 104, 2),
    @('2025-10-13', 'zone tampon', 104, 3),
    @('2025-10-13', 'zone tampon', 105, 2),
    @('2025-10-13', 'agriculture biologique', 106, 2),
    @('2025-10-13', 'ruissellement', 107, 1),
    @('2025-10-13', 'zone tampon', 107, 1),
    @('2025-10-13', 'zone tampon', 108, 1),
    @('2025-10-13', 'zone tampon', 109, 3),
    @('2025-10-13', 'zone tampon', 110, 1),
    @('2025-10-13', 'zone tampon', 112, 1),
    @('2025-10-13', 'eaux souterraines', 114, 2),
    @('2025-10-13', 'zone tampon', 114, 1),
    @('2025-10-13', 'zone tampon', 115, 1),
    @('2025-10-13', 'zone tampon', 116, 1),
    @('2025-10-13', 'zone tampon', 118, 1),
    @('2025-10-13', 'eaux souterraines', 122, 1),
    @('2025-10-13', 'eaux souterraines', 126, 1),
    @('2025-10-13', 'zone tampon', 127, 1),
    @('2025-10-13', 'eaux souterraines', 130, 4),
    @('2025-10-13', 'ruissellement', 131, 2),
    @('2025-10-13', 'eaux souterraines', 131, 1),
    @('2025-10-13', 'zone tampon', 132, 1),
    @('2025-10-13', 'zone tampon', 133, 1),
    @('2025-10-13', 'ruissellement', 188, 1),
    @('2025-10-13', 'ruissellement', 190, 1),
    @('2025-10-13', 'eaux souterraines', 191, 7),
    @('2025-10-13', 'eaux de surface', 191, 1),
    @('2025-10-13', 'eaux souterraines', 194, 2),
    @('2025-10-13', 'eaux souterraines', 195, 4),
    @('2025-10-13', 'eaux de surface', 195, 1),
    @('2025-10-13', 'eaux souterraines', 198, 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Column A holds a date formatted as plain text (e.g. "2025-10-13"),
    # not an Excel date serial. Force text number-format before assigning
    # so Excel does not auto-convert the string to a date, then clear the
    # formatting again so the cell keeps no explicit style (matching the
    # rest of the sheet).
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
